$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (C) column dates from 2023-09-15 (45184) to 2023-09-16 (45185)
# for rows 2 through 8.
$ws.Range("C2:C8").Value = 45185
